$wb = $excel.ActiveWorkbook

# --- Sheet "line_imp": replace literal values with formulas that halve them ---
$ws = $wb.Worksheets.Item("line_imp")

$ws.Range("E2").Formula = "=0.05/2"
$ws.Range("B3").Formula = "=0.05/2"
$ws.Range("C3").Formula = "=0.25/2"
$ws.Range("E3").Formula = "=0.05/2/2"
$ws.Range("E4").Formula = "=0.033333/2"
$ws.Range("E5").Formula = "=0.033333/2"
$ws.Range("E6").Formula = "=0.02/2"

# --- Update selections (view state) on both sheets ---
$ws1 = $wb.Worksheets.Item("initial")
$ws1.Activate() | Out-Null
$ws1.Range("E2").Select() | Out-Null

$ws2 = $wb.Worksheets.Item("line_imp")
$ws2.Activate() | Out-Null
$ws2.Range("E4").Select() | Out-Null
